$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------
# Insert the new "2022-Q1" sheet right after "2021-Q4".
#
# We build it by duplicating the existing "总计" sheet rather than
# Worksheets.Add() -- a straight Add() comes back as a bare-bones sheet
# (no <sheetPr>/outline info, different page margins) whereas Copy()
# gives us a faithful clone (sheetPr, margins, the bold/bordered header
# style already sitting at style index used by 总计) that we then grow
# from 4 columns to 8 and overwrite with the 2022-Q1 figures.
# --------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$total.Copy($null, $q4)

# NOTE: worksheet handles resolved *before* a sheet is inserted/copied
# track sheet *position*, not identity -- so both sheets must be looked
# up again now that "总计" has shifted one slot to the right.
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q1"
$total = $wb.Worksheets.Item("总计")

# Extend the cloned header/index styling from 4 columns (A:D) to 8 (A:H)
$newSheet.Range("D1").Copy($newSheet.Range("E1:H1"))
$newSheet.Range("A2").Copy($newSheet.Range("A3"))

# --- Header row -----------------------------------------------------------
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# --- Row 2 : 001303 / 银华稳利灵活配置混合A --------------------------------
$newSheet.Cells.Item(2,1).Value = 0

$newSheet.Cells.Item(2,2).NumberFormat = "@"
$newSheet.Cells.Item(2,2).Value = "001303"
$newSheet.Cells.Item(2,3).Value = "银华稳利灵活配置混合A"

$newSheet.Cells.Item(2,4).NumberFormat = "@"
$newSheet.Cells.Item(2,4).Value = "0.18"
$newSheet.Cells.Item(2,5).NumberFormat = "@"
$newSheet.Cells.Item(2,5).Value = "28.88"
$newSheet.Cells.Item(2,6).NumberFormat = "@"
$newSheet.Cells.Item(2,6).Value = "0.65"
$newSheet.Cells.Item(2,7).NumberFormat = "@"
$newSheet.Cells.Item(2,7).Value = "0.0012"

$newSheet.Cells.Item(2,8).Value = 6

# --- Row 3 : 002323 / 银华稳利灵活配置混合C --------------------------------
$newSheet.Cells.Item(3,1).Value = 1

$newSheet.Cells.Item(3,2).NumberFormat = "@"
$newSheet.Cells.Item(3,2).Value = "002323"
$newSheet.Cells.Item(3,3).Value = "银华稳利灵活配置混合C"

$newSheet.Cells.Item(3,4).NumberFormat = "@"
$newSheet.Cells.Item(3,4).Value = "0.12"
$newSheet.Cells.Item(3,5).NumberFormat = "@"
$newSheet.Cells.Item(3,5).Value = "28.88"
$newSheet.Cells.Item(3,6).NumberFormat = "@"
$newSheet.Cells.Item(3,6).Value = "0.65"
$newSheet.Cells.Item(3,7).NumberFormat = "@"
$newSheet.Cells.Item(3,7).Value = "0.0008"

$newSheet.Cells.Item(3,8).Value = 6

# --------------------------------------------------------------------
# Update the "总计" (totals) sheet: push the old 2021-Q4 row down to row
# 3 (carrying its style along), then write the new 2022-Q1 summary row
# into row 2.
# --------------------------------------------------------------------
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2021-Q4"
$total.Cells.Item(3,3).Value = 1
$total.Cells.Item(3,4).Value = 0

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0

# Leave the original sheet active/selected, as it was before the edit.
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Activate()
